# Add ability-ID list columns to the Individuals sheet, and load values
# for each individual row, per the "load abilities to individuals" commit.

$wb = $excel.ActiveWorkbook
$wsIndividuals = $wb.Worksheets.Item("Individuals")
$wsTargeted = $wb.Worksheets.Item("TargetedAbilities")

# --- Insert 4 new columns (AT:AW) in front of the old AT column ---------
# (this shifts the old AT:BD "animation" columns right by four, to AX:BH)
$wsIndividuals.Columns("AT:AW").Insert()

# --- New column headers ---------------------------------------------------
$wsIndividuals.Range("AT1").Value2 = "permenantAbilityIDs"
$wsIndividuals.Range("AU1").Value2 = "durationAbilityIDs"
$wsIndividuals.Range("AV1").Value2 = "targetedAbilityIDs"
$wsIndividuals.Range("AW1").Value2 = "InstantAbilityIDs"

# --- AU2:AX9 hold ability-id CSV lists, so force Text format on that whole
#     block before typing into it (matches the numFmt "@" style already used
#     elsewhere in this workbook) so comma lists like "2,100,101" aren't
#     reinterpreted as numbers.
$wsIndividuals.Range("AU2:AX9").NumberFormat = "@"

# --- permenantAbilityIDs (AT) - no abilities loaded for any individual yet
$wsIndividuals.Range("AT2:AT9").Value2 = -1

# --- durationAbilityIDs (AU) - none loaded yet
$wsIndividuals.Range("AU2:AU9").Value2 = -1

# --- targetedAbilityIDs (AV) - "2,100,101" loaded for rows 2, 8 and 9
$wsIndividuals.Range("AV2").Value2 = "2,100,101"
$wsIndividuals.Range("AV3:AV7").Value2 = -1
$wsIndividuals.Range("AV8").Value2 = "2,100,101"
$wsIndividuals.Range("AV9").Value2 = "2,100,101"

# --- InstantAbilityIDs (AW) - none loaded yet
$wsIndividuals.Range("AW2:AW9").Value2 = -1

# --- Column widths for the newly-inserted columns --------------------------
$wsIndividuals.Range("AT1").ColumnWidth = 19.7109375
$wsIndividuals.Range("AU1:AV1").ColumnWidth = 17.42578125
$wsIndividuals.Range("AW1").ColumnWidth = 15.85546875

# --- View/selection bookkeeping -------------------------------------------
# Update the (previously active) TargetedAbilities sheet's selection first …
$wsTargeted.Range("A2:P3").Select()

# … then make Individuals the active sheet/tab again, with the new columns
# scrolled into view, matching where the edit was made.
$wsIndividuals.Activate()
$wsIndividuals.Range("AW6").Select()
$excel.ActiveWindow.ScrollColumn = 47
$excel.ActiveWindow.ScrollRow = 1
